$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update gene (column) headers: S100A8 -> CD14, GZMA -> NKG7
$ws.Range("C1").Value = "CD14"
$ws.Range("D1").Value = "NKG7"

# Update count values that changed between the two uploaded versions
$ws.Range("D2").Value = 1
$ws.Range("C3").Value = 18
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 160

# The re-uploaded sheet has no explicit cell formatting (the extra cell
# style used by the old file is gone) and A1 is completely empty rather
# than holding an empty string - replicate that.
$ws.Cells.Style = "Normal"
$ws.Range("A1").Clear()
